$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Removing the 'any' identifier type: its row (A4, "any") goes away and every
# type label below it moves up one slot in the A column (A5->A4, A6->A5, ...,
# A10->A9), with a new "reserved" label taking the vacated A10 slot. The
# numeric bit-flag grid (B:J) underneath is a positional bit-counter that
# does not depend on the label, so it is left untouched.
$ws.Range("A4").Value = "string"
$ws.Range("A5").Value = "boolean"
$ws.Range("A6").Value = "integer"
$ws.Range("A7").Value = "float"
$ws.Range("A8").Value = "long"
$ws.Range("A9").Value = "bytes"
$ws.Range("A10").Value = "reserved"

# The SEMANTIC sub-table (L:P, rows 4-6) references the bit-flag grid by
# cell address (uuid = fixed-list + <row that used to hold "long">, etc.).
# Since "long" (and everything after it) moved up one row, those references
# shift from J9/J20 to J8/J19-relative addresses.
$ws.Range("O4").Formula = "=J21+J8"
$ws.Range("O5").Formula = "=J8"
$ws.Range("O6").Formula = "=J20+J7"

# Selection moves as recorded in the saved view state.
$ws.Range("O7").Select()

$wb.Save()
